$wb = $excel.ActiveWorkbook

# --- Add 1 new Startup (SU) company to the "Startup" sheet, row 36 ---
$wsStartup = $wb.Worksheets.Item("Startup")
$wsStartup.Activate()

$wsStartup.Range("B36").Value = "Arago"
$wsStartup.Range("C36").Value = 15
$wsStartup.Range("D36").Value = "Si Photonics for AI acceleration (in stealth)"
$wsStartup.Range("E36").Value = "Hans-Christian Boos"
$wsStartup.Range("F36").Value = "Paris, France and San Francisco, California"

$wsStartup.Range("F36").Select()

# --- Add 1 new Defense Contractor (DC) to the "Defense Contractors" sheet, row 6 ---
$wsDefense = $wb.Worksheets.Item("Defense Contractors")
$wsDefense.Activate()

$wsDefense.Range("B6").Value = "Critical Frequency Design"
$wsDefense.Range("C6").Value = 40
$wsDefense.Range("D6").Value = "Microwave photonics for communications and sensing. Frequency Specific Limiter via SBS. Free Space Optics"
$wsDefense.Range("E6").Value = "Johnathen Warren and Charles Middleton "
$wsDefense.Range("F6").Value = "Melbourne, FL, USA"

$wsDefense.Range("E7").Select()
